$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
$ws.Range("A8").Value = "Volume 29   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/14/2022  Through  11/20/2022"

# --- Data table updates (rows 14-27) --------------------------------------
$ws.Range("N14").Value = -80.95238095238
$ws.Range("D15").Value = 1
$ws.Range("F15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("K15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1
$ws.Range("F15").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").Value = 0
$ws.Range("K15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = 58.333333333333
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -85.714285714285
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 22.222222222222
$ws.Range("I16").Value = 202
$ws.Range("J16").Value = 173
$ws.Range("K16").Value = 16.763005780346
$ws.Range("L16").Value = 28.662420382165
$ws.Range("M16").Value = -14.767932489451
$ws.Range("N16").Value = -74.135723431498
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -62.5
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = -48.648648648648
$ws.Range("I17").Value = 309
$ws.Range("J17").Value = 312
$ws.Range("K17").Value = -0.961538461538
$ws.Range("L17").Value = 20.703125
$ws.Range("M17").Value = 134.090909090909
$ws.Range("N17").Value = -18.037135278514
$ws.Range("C18").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -62.5
$ws.Range("J18").Value = 97
$ws.Range("K18").Value = 38.144329896907
$ws.Range("L18").Value = -7.586206896551
$ws.Range("M18").Value = -49.812734082397
$ws.Range("N18").Value = -87.417840375586
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -23.076923076923
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = -3.703703703703
$ws.Range("I19").Value = 586
$ws.Range("J19").Value = 453
$ws.Range("K19").Value = 29.359823399558
$ws.Range("L19").Value = 22.594142259414
$ws.Range("M19").Value = 83.125
$ws.Range("N19").Value = 10.566037735849
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -16.666666666666
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = -4.347826086956
$ws.Range("I20").Value = 200
$ws.Range("J20").Value = 176
$ws.Range("K20").Value = 13.636363636363
$ws.Range("L20").Value = 21.212121212121
$ws.Range("M20").Value = -26.739926739926
$ws.Range("N20").Value = -93.160054719562
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = -51.282051282051
$ws.Range("F21").Value = 122
$ws.Range("G21").Value = 149
$ws.Range("H21").Value = -18.120805369127
$ws.Range("I21").Value = 1454
$ws.Range("J21").Value = 1228
$ws.Range("K21").Value = 18.403908794788
$ws.Range("L21").Value = 19.868095630667
$ws.Range("M21").Value = 16.227018385291
$ws.Range("N21").Value = -74.611489436004
$ws.Range("C22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 25
$ws.Range("K22").Value = 78.571428571428
$ws.Range("L22").Value = 56.25
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -22.727272727272
$ws.Range("F24").Value = 89
$ws.Range("H24").Value = -16.822429906542
$ws.Range("I24").Value = 1220
$ws.Range("J24").Value = 847
$ws.Range("K24").Value = 44.037780401416
$ws.Range("L24").Value = 59.685863874345
$ws.Range("M24").Value = 98.051948051948
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 42.857142857142
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = 23.076923076923
$ws.Range("I25").Value = 477
$ws.Range("J25").Value = 419
$ws.Range("K25").Value = 13.842482100238
$ws.Range("L25").Value = 20.454545454545
$ws.Range("M25").Value = 6.236080178173
$ws.Range("D26").Value = 1
$ws.Range("F26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = 0
$ws.Range("K26").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 1
$ws.Range("F26").Copy()
$ws.Range("G26").PasteSpecial(-4122)
$ws.Range("H26").Value = 200
$ws.Range("K26").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("I26").Value = 35
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = 94.444444444444
$ws.Range("L26").Value = 94.444444444444
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("J27").Value = 36
$ws.Range("K27").Value = 44.444444444444

$excel.CutCopyMode = 0
